$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 19: merge of mice and extratrees_40 ----
$ws.Range("A19").Value = "merge of mice and extratrees_40"
$ws.Range("B19").Value = 0.18781890000000001
$ws.Range("C19").Value = 0.24322659999999999
$ws.Range("D19").Value = 0.21782019999999999
$ws.Range("E19").Value = 0.20580490000000001
$ws.Range("F19").Value = 0.1467821
$ws.Range("G19").Value = 0.14559250000000001
$ws.Range("H19").Value = 0.25045699999999999
$ws.Range("I19").Value = 0.196493
$ws.Range("J19").Value = 0.22881180000000001
$ws.Range("K19").Value = 0.2147213
$ws.Range("L19").Value = 0.18276020000000001
$ws.Range("M19").Value = 0.2105648
$ws.Range("N19").Value = 0.26546540000000002

# ---- Row 20: 3D-mice (iter=2, seed=100) ----
$ws.Range("A20").Value = "3D-mice (iter=2, seed=100)"
$ws.Range("B20").Value = 0.2028953
$ws.Range("C20").Value = 0.2625016
$ws.Range("D20").Value = 0.2355701
$ws.Range("E20").Value = 0.21756149999999999
$ws.Range("F20").Value = 0.14724799999999999
$ws.Range("G20").Value = 0.1457425
$ws.Range("H20").Value = 0.27084839999999999
$ws.Range("I20").Value = 0.2311107
$ws.Range("J20").Value = 0.25933030000000001
$ws.Range("K20").Value = 0.24842620000000001
$ws.Range("L20").Value = 0.1901919
$ws.Range("M20").Value = 0.2330045
$ws.Range("N20").Value = 0.2771246

# ---- Row 21: 3D-mice (iter=1, seed=100) ----
$ws.Range("A21").Value = "3D-mice (iter=1, seed=100)"
$ws.Range("B21").Value = 0.2024956
$ws.Range("C21").Value = 0.26413900000000001
$ws.Range("D21").Value = 0.23642199999999999
$ws.Range("E21").Value = 0.217616
$ws.Range("F21").Value = 0.1467599
$ws.Range("G21").Value = 0.14528869999999999
$ws.Range("H21").Value = 0.27149210000000001
$ws.Range("I21").Value = 0.23414409999999999
$ws.Range("J21").Value = 0.26195269999999998
$ws.Range("K21").Value = 0.2517971
$ws.Range("L21").Value = 0.19084180000000001
$ws.Range("M21").Value = 0.23372229999999999
$ws.Range("N21").Value = 0.27795180000000003

# ---- Fix up cell formats so styles match the recorded layout exactly ----
# Source cells that already carry the three distinct formats used across B:N.
$fmtPlain  = $ws.Range("L2")   # s=4 (black Arial number format)
$fmtTheme  = $ws.Range("B1")   # s=5 (theme-color Arial number format)
$fmtLeft   = $ws.Range("N16")  # s=8 (left/vcenter, no number format)

$fmtPlain.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("I20").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("I21").PasteSpecial(-4122)

$fmtTheme.Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("E21").PasteSpecial(-4122)

$fmtLeft.Copy()
$ws.Range("B20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---- Column B is now a touch wider than the rest of the data columns ----
$ws.Columns("B").ColumnWidth = 10.5

# ---- Active cell moves to D9 ----
$ws.Range("D9").Select()
